# Applies the crypto-list refresh from the GitHub Actions update:
# new prices/volume percentages for existing coins, plus a ranking
# reshuffle among rows 48-51 (WEMIXToken / Maker / Stellar / ApeXProtocol).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel auto-converts plain decimal-looking text ("1.00", "11.00", ...) into
# a Number, which would silently drop the trailing zeros the site displays.
# Force those cells to stay Text the same way a user would (leading quote),
# then restore the default "Normal" style so we do not leave a stray
# quote-prefixed style behind on the cell.
function Set-TextCell($addr, $value) {
    $ws.Range($addr).Value = "'" + $value
    $ws.Range($addr).Style = "Normal"
}

$ws.Range("D2").Value = "67.942.92"
$ws.Range("E2").Value = "  -3.60%  "

$ws.Range("D3").Value = "3.685.52"
$ws.Range("E3").Value = "  -4.32%  "

$ws.Range("E4").Value = "  +0.12%  "

Set-TextCell "D5" "595.79"
$ws.Range("E5").Value = "  +1.68%  "

Set-TextCell "D6" "180.17"
$ws.Range("E6").Value = "  +8.64%  "

$ws.Range("D7").Value = "3.679.06"
$ws.Range("E7").Value = "  -4.24%  "

Set-TextCell "D8" "0.624"
$ws.Range("E8").Value = "  -6.24%  "

Set-TextCell "D9" "1.00"
$ws.Range("E9").Value = "  +0.11%  "

Set-TextCell "D10" "0.715"
$ws.Range("E10").Value = "  -3.95%  "

Set-TextCell "D11" "0.162"
$ws.Range("E11").Value = "  -6.44%  "

Set-TextCell "D12" "55.90"
$ws.Range("E12").Value = "  +5.27%  "

Set-TextCell "D13" "0.0000288"
$ws.Range("E13").Value = "  -9.79%  "

Set-TextCell "D14" "10.35"
$ws.Range("E14").Value = "  -7.56%  "

$ws.Range("D15").Value = "4.254.31"
$ws.Range("E15").Value = "  -4.62%  "

$ws.Range("D16").Value = "3.683.03"
$ws.Range("E16").Value = "  -4.87%  "

Set-TextCell "D17" "19.34"
$ws.Range("E17").Value = "  -6.16%  "

$ws.Range("E18").Value = "  -2.27%  "

Set-TextCell "D19" "12.78"
$ws.Range("E19").Value = "  -7.13%  "

$ws.Range("E20").Value = "  -6.56%  "

$ws.Range("D21").Value = "67.729.24"
$ws.Range("E21").Value = "  -3.60%  "

Set-TextCell "D22" "407.74"
$ws.Range("E22").Value = "  -6.14%  "

Set-TextCell "D23" "4.53"
$ws.Range("E23").Value = "  -3.09%  "

Set-TextCell "D24" "88.18"
$ws.Range("E24").Value = "  -5.95%  "

$ws.Range("E25").Value = "  -8.03%  "

Set-TextCell "D26" "12.75"
$ws.Range("E26").Value = "  -6.85%  "

Set-TextCell "D27" "11.00"
$ws.Range("E27").Value = "  +1.14%  "

Set-TextCell "D28" "4.00"
$ws.Range("E28").Value = "  -1.37%  "

Set-TextCell "D29" "6.07"
$ws.Range("E29").Value = "  +2.32%  "

Set-TextCell "D30" "9.41"
$ws.Range("E30").Value = "  -7.85%  "

Set-TextCell "D31" "32.58"
$ws.Range("E31").Value = "  -6.64%  "

$ws.Range("E32").Value = "  -8.19%  "

Set-TextCell "D33" "12.44"
$ws.Range("E33").Value = "  -7.34%  "

$ws.Range("E34").Value = "  -6.30%  "

Set-TextCell "D35" "43.30"
$ws.Range("E35").Value = "  -10.73%  "

Set-TextCell "D36" "64.02"
$ws.Range("E36").Value = "  -7.66%  "

Set-TextCell "D37" "589.86"
$ws.Range("E37").Value = "  -4.06%  "

$ws.Range("D38").Value = "0.0₃0871"
$ws.Range("E38").Value = "  -10.06%  "

$ws.Range("E39").Value = "  +0.06%  "

Set-TextCell "D40" "0.399"
$ws.Range("E40").Value = "  -4.19%  "

Set-TextCell "D41" "0.997"
$ws.Range("E41").Value = "  -0.22%  "

Set-TextCell "D42" "0.136"
$ws.Range("E42").Value = "  -4.50%  "

Set-TextCell "D43" "2.76"
$ws.Range("E43").Value = "  +3.51%  "

$ws.Range("E44").Value = "  -7.94%  "

Set-TextCell "D45" "0.0432"
$ws.Range("E45").Value = "  -6.75%  "

Set-TextCell "D46" "2.88"
$ws.Range("E46").Value = "  -10.36%  "

Set-TextCell "D47" "9.18"
$ws.Range("E47").Value = "  -7.87%  "

$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell "D48" "2.71"
$ws.Range("E48").Value = "  -2.55%  "

$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.763.88"
$ws.Range("E49").Value = "  -2.37%  "

$ws.Range("E50").Value = "  -6.46%  "

$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextCell "D51" "3.13"
$ws.Range("E51").Value = "  -4.63%  "
